$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '89.314.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +10.09%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.371.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +7.31%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '649.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.43%  '

$ws.Range("E7").Value = '  +47.32%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.616'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.366.71'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.617'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000293'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +16.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +15.88%  '

$ws.Range("E14").Value = '  +2.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.981.05'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.05%  '

$ws.Range("E16").Value = '  +5.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.182.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +10.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.366.73'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.10%  '

$ws.Range("E19").Value = '  +7.25%  '

$ws.Range("E20").Value = '  -0.98%  '

$ws.Range("E21").Value = '  +8.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '459.09'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +18.41%  '

$ws.Range("E27").Value = '  +6.77%  '

$ws.Range("E28").Value = '  +19.59%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '79.52'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.198'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +44.03%  '

$ws.Range("E31").Value = '  -0.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.39'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '598.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.75%  '

$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.59'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.93%  '

$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.21%  '

$ws.Range("E36").Value = '  +7.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +20.26%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.144'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.86%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.430'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.45%  '

$ws.Range("E41").Value = '  +6.20%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.93%  '

$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.86'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.10%  '

$ws.Range("E45").Value = '  +11.18%  '

$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.65'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '188.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '46.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.14%  '

$ws.Range("E50").Value = '  +8.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.669'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.45%  '
